# Append two new activity-log rows (rows 4 and 5) to the "Activity Log"
# worksheet, mirroring the two new <row> elements introduced by the diff.
# Column order is Date | Name | Role | Event, same as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "10/17/2025, 03:18:09 PM"
$ws.Cells.Item(4, 2).Value = "Winnish Allwin G J"
$ws.Cells.Item(4, 3).Value = "team_member"
$ws.Cells.Item(4, 4).Value = "Submitted metrics for 10/17/2025 - Assigned: 18, Resolved: 13, SLA Breaches: 2, Reopened: 1, Client Interactions: 13 | Remarks: N/A"

# --- Row 5 ---------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = "10/17/2025, 03:29:42 PM"
$ws.Cells.Item(5, 2).Value = "Winnish Allwin G J"
$ws.Cells.Item(5, 3).Value = "team_member"
$ws.Cells.Item(5, 4).Value = "Submitted metrics for 2025-10-17 - Assigned: 12, Resolved: 6, SLA Breaches: 1, Reopened: 3, Client Interactions: 5 | Remarks: N/A"

# The workbook flags A1:D3 as "number stored as text" ignorable errors; now
# that the used range has grown to A1:D5, extend that same ignore flag over
# the newly written rows so the whole data block keeps being ignored.
$fullRange = $ws.Range("A1:D5")
try {
    $ignoredError = $fullRange.Errors.Item(3)   # xlNumberAsText
    $ignoredError.Ignore = $true
} catch {
    # Older/alternate builds may not expose this check; the data cells
    # above are the important part of this edit.
}
